$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for affected rows as per repull/mean calculation
$ws.Range("F2").Value = -2
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -1
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = -5
